$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E10").Value = 25
$ws.Range("E14").Value = 34
$ws.Range("E15").Value = 88
$ws.Range("E16").Value = 298
$ws.Range("E18").Value = 89
